$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

$ws.Range("H125").Value = 2850
$ws.Range("I125").Value = 1575
$ws.Range("K125").Value = 14175
$ws.Range("M125").Value = -11715

$ws.Range("H129").Value = 3572934.2
$ws.Range("I129").Value = 83335416
$ws.Range("J129").Value = 1479.8806
$ws.Range("K129").Value = 250006248
$ws.Range("L129").Value = 4439.641799999999
$ws.Range("M129").Value = -250001248
$ws.Range("N129").Value = -14439.6418

$ws.Range("H138").Value = 3362.7937
$ws.Range("I138").Value = 1458.6097
$ws.Range("J138").Value = 6911.5
$ws.Range("K138").Value = 4375.8291
$ws.Range("L138").Value = 20734.5
$ws.Range("M138").Value = 764.1709000000001
$ws.Range("N138").Value = -31014.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4642.8096
$ws.Range("I61").Value = 1722.6666
$ws.Range("J61").Value = 6832.9165
$ws.Range("K61").Value = 1722.6666
$ws.Range("L61").Value = 6832.9165
$ws.Range("M61").Value = -1510.6666
$ws.Range("N61").Value = -7256.9165

$ws.Range("H122").Value = 2960.75
$ws.Range("I122").Value = 2020.4117
$ws.Range("K122").Value = 6061.2351
$ws.Range("M122").Value = -3611.2351

$ws.Range("H132").Value = 19611098
$ws.Range("I132").Value = 25644104
$ws.Range("K132").Value = 76932312
$ws.Range("M132").Value = -76929782

$ws.Range("H136").Value = 4642.8096
$ws.Range("I136").Value = 1722.6666
$ws.Range("J136").Value = 6832.9165
$ws.Range("K136").Value = 5167.9998
$ws.Range("L136").Value = 20498.7495
$ws.Range("M136").Value = -2617.9998
$ws.Range("N136").Value = -25598.7495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1528.9375
$ws.Range("I107").Value = 450
$ws.Range("K107").Value = 450
$ws.Range("M107").Value = 1470

$ws.Range("H134").Value = 3036.9355
$ws.Range("I134").Value = 1832.5652
$ws.Range("J134").Value = 6499.5
$ws.Range("K134").Value = 5497.6956
$ws.Range("L134").Value = 19498.5
$ws.Range("M134").Value = -2962.6956
$ws.Range("N134").Value = -24568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 581500.4399999999
$ws.Range("I6").Value = 1334500.4
$ws.Range("J6").Value = 16750.5
$ws.Range("K6").Value = 1334500.4
$ws.Range("L6").Value = 16750.5
$ws.Range("M6").Value = -1334387.4
$ws.Range("N6").Value = -16976.5

$ws.Range("H22").Value = 1658.875
$ws.Range("I22").Value = 194.2
$ws.Range("K22").Value = 194.2
$ws.Range("M22").Value = 155.8

$ws.Range("H26").Value = 70021
$ws.Range("J26").Value = 70021
$ws.Range("L26").Value = 70021
$ws.Range("N26").Value = -70595

$ws.Range("H31").Value = 2630.9048
$ws.Range("I31").Value = 1600.9166
$ws.Range("J31").Value = 4004.2222
$ws.Range("K31").Value = 1600.9166
$ws.Range("L31").Value = 4004.2222
$ws.Range("M31").Value = -1305.9166
$ws.Range("N31").Value = -4594.2222

$ws.Range("H34").Value = 2630.9048
$ws.Range("I34").Value = 1600.9166
$ws.Range("J34").Value = 4004.2222
$ws.Range("K34").Value = 1600.9166
$ws.Range("L34").Value = 4004.2222
$ws.Range("M34").Value = -1398.9166
$ws.Range("N34").Value = -4408.2222

$ws.Range("H107").Value = 2659.7144
$ws.Range("I107").Value = 644.4286
$ws.Range("J107").Value = 4675
$ws.Range("K107").Value = 644.4286
$ws.Range("L107").Value = 4675
$ws.Range("M107").Value = 1275.5714
$ws.Range("N107").Value = -8515

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1084.0975
$ws.Range("I107").Value = 666.8461
$ws.Range("J107").Value = 1277.8214
$ws.Range("K107").Value = 2000.5383
$ws.Range("L107").Value = 3833.4642
$ws.Range("M107").Value = -80.53829999999994
$ws.Range("N107").Value = -7673.4642

$ws.Range("H120").Value = 14116.25
$ws.Range("I120").Value = 8232.5
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 24697.5
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -19859.5
$ws.Range("N120").Value = -69676

$ws.Range("H129").Value = 38695.855
$ws.Range("I129").Value = 5915.75
$ws.Range("J129").Value = 51807.9
$ws.Range("K129").Value = 17747.25
$ws.Range("L129").Value = 155423.7
$ws.Range("M129").Value = -12747.25
$ws.Range("N129").Value = -165423.7

$ws.Range("H131").Value = 1450.9231
$ws.Range("J131").Value = 1373.2
$ws.Range("L131").Value = 4119.6
$ws.Range("N131").Value = -14199.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2036.4231
$ws.Range("I102").Value = 1284.1875
$ws.Range("J102").Value = 3240
$ws.Range("K102").Value = 1284.1875
$ws.Range("L102").Value = 3240
$ws.Range("M102").Value = 337.8125
$ws.Range("N102").Value = -6484

$ws.Range("H113").Value = 4587.2856
$ws.Range("I113").Value = 2455.5
$ws.Range("J113").Value = 5440
$ws.Range("K113").Value = 2455.5
$ws.Range("L113").Value = 5440
$ws.Range("M113").Value = -285.5
$ws.Range("N113").Value = -9780

$ws.Range("H132").Value = 33336816
$ws.Range("I132").Value = 62502156
$ws.Range("K132").Value = 187506468
$ws.Range("M132").Value = -187503938

$ws.Range("H133").Value = 28426.666
$ws.Range("J133").Value = 28426.666
$ws.Range("L133").Value = 28426.666
$ws.Range("N133").Value = -38546.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2957.1428
$ws.Range("I122").Value = 2568.4211
$ws.Range("J122").Value = 3777.7778
$ws.Range("K122").Value = 7705.263300000001
$ws.Range("L122").Value = 11333.3334
$ws.Range("M122").Value = -5255.263300000001
$ws.Range("N122").Value = -16233.3334

$ws.Range("H132").Value = 3662.15
$ws.Range("I132").Value = 3158.3333
$ws.Range("J132").Value = 3878.0715
$ws.Range("K132").Value = 9474.999899999999
$ws.Range("L132").Value = 11634.2145
$ws.Range("M132").Value = -6944.999899999999
$ws.Range("N132").Value = -16694.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1734.3636
$ws.Range("I113").Value = 1192.25
$ws.Range("J113").Value = 2044.1428
$ws.Range("K113").Value = 3576.75
$ws.Range("L113").Value = 6132.428400000001
$ws.Range("M113").Value = -1406.75
$ws.Range("N113").Value = -10472.4284

$ws.Range("H122").Value = 2576
$ws.Range("I122").Value = 1953.909
$ws.Range("K122").Value = 5861.727000000001
$ws.Range("M122").Value = -3411.727000000001
